$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 2933294.8
$ws.Range("I33").Value = 791.3103599999999
$ws.Range("J33").Value = 45454590
$ws.Range("K33").Value = 791.3103599999999
$ws.Range("L33").Value = 45454590
$ws.Range("M33").Value = -562.3103599999999
$ws.Range("N33").Value = -45455048

# Row 69
$ws.Range("H69").Value = 4482.4116
$ws.Range("I69").Value = 4023.2
$ws.Range("J69").Value = 4673.75
$ws.Range("K69").Value = 12069.6
$ws.Range("L69").Value = 14021.25
$ws.Range("M69").Value = -11195.6
$ws.Range("N69").Value = -15769.25

# Row 72
$ws.Range("H72").Value = 4482.4116
$ws.Range("I72").Value = 4023.2
$ws.Range("J72").Value = 4673.75
$ws.Range("K72").Value = 36208.8
$ws.Range("L72").Value = 42063.75
$ws.Range("M72").Value = -31840.8
$ws.Range("N72").Value = -50799.75

# Row 123
$ws.Range("H123").Value = 28831.7
$ws.Range("J123").Value = 28831.7
$ws.Range("L123").Value = 28831.7
$ws.Range("N123").Value = -38631.7

# Row 130
$ws.Range("H130").Value = 43375
$ws.Range("J130").Value = 43375
$ws.Range("L130").Value = 43375
$ws.Range("N130").Value = -53415

# Row 132
$ws.Range("H132").Value = 1247.0625
$ws.Range("I132").Value = 1130.2
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 3390.6
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -860.6000000000004
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3990.1948
$ws.Range("I32").Value = 2962.1016
$ws.Range("K32").Value = 2962.1016
$ws.Range("M32").Value = -2675.1016

# Row 45
$ws.Range("H45").Value = 17885.334
$ws.Range("I45").Value = 17885.334
$ws.Range("K45").Value = 17885.334
$ws.Range("M45").Value = -17508.334

# Row 74
$ws.Range("H74").Value = 1299.2142
$ws.Range("I74").Value = 1214.6316
$ws.Range("J74").Value = 1477.7778
$ws.Range("K74").Value = 1214.6316
$ws.Range("L74").Value = 1477.7778
$ws.Range("M74").Value = -340.6315999999999
$ws.Range("N74").Value = -3225.7778

# Row 77
$ws.Range("H77").Value = 1299.2142
$ws.Range("I77").Value = 1214.6316
$ws.Range("J77").Value = 1477.7778
$ws.Range("K77").Value = 6073.157999999999
$ws.Range("L77").Value = 7388.889
$ws.Range("M77").Value = -1705.157999999999
$ws.Range("N77").Value = -16124.889

$ws = $wb.Worksheets.Item("BSM")
# Row 74
$ws.Range("H74").Value = 20261
$ws.Range("J74").Value = 23794
$ws.Range("L74").Value = 23794
$ws.Range("N74").Value = -25666

# Row 77
$ws.Range("H77").Value = 20261
$ws.Range("J77").Value = 23794
$ws.Range("L77").Value = 71382
$ws.Range("N77").Value = -80742

# Row 81
$ws.Range("H81").Value = 34916.668
$ws.Range("J81").Value = 34916.668
$ws.Range("L81").Value = 34916.668
$ws.Range("N81").Value = -37038.668

# Row 84
$ws.Range("H84").Value = 34916.668
$ws.Range("J84").Value = 34916.668
$ws.Range("L84").Value = 104750.004
$ws.Range("N84").Value = -115358.004

# Row 94
$ws.Range("H94").Value = 1255.3871
$ws.Range("I94").Value = 876.6799999999999
$ws.Range("J94").Value = 2833.3333
$ws.Range("K94").Value = 876.6799999999999
$ws.Range("L94").Value = 2833.3333
$ws.Range("M94").Value = -425.6799999999999
$ws.Range("N94").Value = -3735.3333

# Row 130
$ws.Range("H130").Value = 45000
$ws.Range("J130").Value = 45000
$ws.Range("L130").Value = 45000
$ws.Range("N130").Value = -55040

# Row 139
$ws.Range("H139").Value = 35280
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 35280
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 35280
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -45560

$ws = $wb.Worksheets.Item("CRP")
# Row 64
$ws.Range("H64").Value = 39466.668
$ws.Range("J64").Value = 39466.668
$ws.Range("L64").Value = 39466.668
$ws.Range("N64").Value = -39962.668

# Row 67
$ws.Range("H67").Value = 39466.668
$ws.Range("J67").Value = 39466.668
$ws.Range("L67").Value = 39466.668
$ws.Range("N67").Value = -41182.668

# Row 98
$ws.Range("H98").Value = 30493.334
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 30493.334
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 30493.334
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -34985.334

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 9000087
$ws.Range("I12").Value = 12500106
$ws.Range("J12").Value = 7353019.5
$ws.Range("K12").Value = 37500318
$ws.Range("L12").Value = 22059058.5
$ws.Range("M12").Value = -37500145
$ws.Range("N12").Value = -22059404.5

# Row 32
$ws.Range("H32").Value = 1441.4286
$ws.Range("J32").Value = 1515
$ws.Range("L32").Value = 4545
$ws.Range("N32").Value = -5111

# Row 38
$ws.Range("H38").Value = 5263184
$ws.Range("I38").Value = 7142875.5
$ws.Range("J38").Value = 48
$ws.Range("K38").Value = 21428626.5
$ws.Range("L38").Value = 144
$ws.Range("M38").Value = -21428279.5
$ws.Range("N38").Value = -838

# Row 54
$ws.Range("H54").Value = 2973.6843
$ws.Range("J54").Value = 2973.6843
$ws.Range("L54").Value = 8921.052899999999
$ws.Range("N54").Value = -10039.0529

# Row 87
$ws.Range("H87").Value = 1899.8
$ws.Range("I87").Value = 1899.8
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 5699.4
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -4451.4
$ws.Range("N87").ClearContents()

# Row 90
$ws.Range("H90").Value = 1899.8
$ws.Range("I90").Value = 1899.8
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 17098.2
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -10858.2
$ws.Range("N90").ClearContents()

# Row 105
$ws.Range("H105").Value = 6728.6113
$ws.Range("J105").Value = 6799.353
$ws.Range("L105").Value = 20398.059
$ws.Range("N105").Value = -25640.059

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 569.8261
$ws.Range("I16").Value = 562.3
$ws.Range("J16").Value = 620
$ws.Range("K16").Value = 562.3
$ws.Range("L16").Value = 620
$ws.Range("M16").Value = -392.3
$ws.Range("N16").Value = -960

# Row 22
$ws.Range("H22").Value = 2779698.5
$ws.Range("I22").Value = 22222718
$ws.Range("J22").Value = 2124.2
$ws.Range("K22").Value = 22222718
$ws.Range("L22").Value = 2124.2
$ws.Range("M22").Value = -22222423
$ws.Range("N22").Value = -2714.2

# Row 27
$ws.Range("H27").Value = 2779698.5
$ws.Range("I27").Value = 22222718
$ws.Range("J27").Value = 2124.2
$ws.Range("K27").Value = 22222718
$ws.Range("L27").Value = 2124.2
$ws.Range("M27").Value = -22222611
$ws.Range("N27").Value = -2338.2

# Row 55
$ws.Range("H55").Value = 75000270
$ws.Range("I55").Value = 90909380
$ws.Range("J55").Value = 55555816
$ws.Range("K55").Value = 90909380
$ws.Range("L55").Value = 55555816
$ws.Range("M55").Value = -90909207
$ws.Range("N55").Value = -55556162

# Row 127
$ws.Range("H127").Value = 42199.4
$ws.Range("J127").Value = 42199.4
$ws.Range("L127").Value = 42199.4
$ws.Range("N127").Value = -52119.4

# Row 132
$ws.Range("H132").Value = 6645808
$ws.Range("I132").Value = 10662955
$ws.Range("J132").Value = 2065.8462
$ws.Range("K132").Value = 31988865
$ws.Range("L132").Value = 6197.5386
$ws.Range("M132").Value = -31986335
$ws.Range("N132").Value = -11257.5386

# Row 136
$ws.Range("H136").Value = 7558.9033
$ws.Range("I136").Value = 8230.611000000001
$ws.Range("J136").Value = 6628.846
$ws.Range("K136").Value = 24691.833
$ws.Range("L136").Value = 19886.538
$ws.Range("M136").Value = -22141.833
